# Removed Test Case Inter-Dependency
#
# - ProductLoanInput!B2 ("shortname") was the hard-coded numeric literal 4202
#   (duplicating ProductLoanInput!B3 / the "description" field, which created
#   an inter-test dependency). Decouple it into its own distinct text value.
# - The product-name title text (shared by B1 on both sheets) changes its
#   suffix from "-PERIODIC" to "-PER1st".
# - Selection / active-sheet bookkeeping is reset: ProductLoanOutput becomes
#   the active tab, and both sheets' selections collapse back to B1.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update the product name title text on both sheets (shared string).
$newTitle = "4202-RBI-EI-DB-DL-REC-FEE+INTEREST-RNI-FFC-SAR-FFROP-DAILY-1-CTRFD-MD-TR-1-ONTIME-PER1st"
$wsInput.Range("B1").Value = $newTitle
$wsOutput.Range("B1").Value = $newTitle

# Decouple the shortname cell from the shared numeric literal -> distinct text.
$wsInput.Range("B2").Value = "420g"

# Reset selections back to B1 on both sheets, and drop the saved scroll
# position on the input sheet.
$wsInput.Activate()
$wsInput.Range("B1").Select() | Out-Null
$wsInput.Application.ActiveWindow.ScrollRow = 1
$wsInput.Application.ActiveWindow.ScrollColumn = 1

# Make ProductLoanOutput the active tab, with its selection on B1.
$wsOutput.Activate()
$wsOutput.Range("B1").Select() | Out-Null
